# Lichtenstein onion BOM: remove the 3 chip-resistor-array parts (RN1, RN2,
# RN3) and replace them with twelve discrete 10k 0805 resistors (R14-R25).
# Also fill in the previously-blank "Value" for IC2 (MCP23017).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the resistor-network rows (RN1, RN2, RN3) ---
# These currently occupy rows 43:45.
$ws.Rows("43:45").Delete()

# --- Insert 12 fresh rows in their place for the new discrete resistors ---
# (this pushes S1, S2, SV1, SV2, SV3, U1 further down, which matches the
# target layout where those parts end up on rows 55-60)
$ws.Rows("43:54").Insert()

$newParts = @("R14","R15","R16","R17","R18","R19","R20","R21","R22","R23","R24","R25")

# Fill in the Part (column A) names first ...
for ($i = 0; $i -lt $newParts.Length; $i++) {
    $row = 43 + $i
    $ws.Cells.Item($row, 1).Value = $newParts[$i]
}

# ... then the rest of the row data (Value/Device/Package/Description/MOUSER)
for ($i = 0; $i -lt $newParts.Length; $i++) {
    $row = 43 + $i
    $ws.Cells.Item($row, 2).Value = "10k"
    $ws.Cells.Item($row, 3).Value = "R-EU_M0805"
    $ws.Cells.Item($row, 4).Value = "M0805"
    $ws.Cells.Item($row, 5).Value = "RESISTOR, European symbol"
    $ws.Cells.Item($row, 6).Value = "652-CR0805FX-1002ELF"
}

# --- Fill in the Value for IC2 (was blank before) ---
$ws.Range("B22").Value = "MCP23017"

# --- Match the saved selection/view state ---
$ws.Range("B23").Select()
